$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that follows the title
#    (Heading1 "Play Dark Vortex Free Slot Game | Yggdrasil Gaming").
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $p.Range.Delete()
        break
    }
}

# 2. Insert a new bold "Play Dark Vortex Free Slot Game | Yggdrasil Gaming"
#    paragraph right before the final (image-prompt) paragraph.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dark Vortex Free Slot Game | Yggdrasil Gaming</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertResult = $newPara.Range.InsertXML($titleXml)

# 3. Replace the old image-generation-prompt text of the final paragraph with
#    the review blurb (keeping its existing italic run formatting).
$oldText = "Create an enticing feature image for Dark Vortex with the following specifications: Style: Cartoony Subject: A happy Maya warrior with glasses should be the main focus of the image. The warrior should be wearing a detailed headpiece and extravagant clothing that make them stand out from the dark and ominous background. They should be smiling and holding a glowing Vortex symbol to add more excitement to the image. Background: The background should reflect the ominous and mysterious atmosphere of the game. It should feature a portal to another dimension, with eerie purple and green hues permeating the scene. The portal should be slightly open, revealing glimpses of the otherworldly realm beyond. Overall Mood: The image should be striking and eye-catching, capturing the attention of potential players. It should convey the thrill and otherworldly allure of the Dark Vortex game, tempting players to dive into the game and uncover its secrets."
$newText = "Read our review of Dark Vortex, a 5-reels and 243-3,125 paylines slot game packed with unique features and scary theme. Play Dark Vortex free today!"
$findResult = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
